$d = $word.ActiveDocument

$replacements = @(
    @("2025-03-14 Friday", "2025-03-15 Saturday"),
    @("22×73=", "69×46="),
    @("14×19=", "25×64="),
    @("13×49=", "83×14="),
    @("93×48=", "30×88="),
    @("23×22=", "94×64="),
    @("59×92=", "63×96="),
    @("53×88=", "23×45="),
    @("13×79=", "72×76="),
    @("63×98=", "87×60="),
    @("13×39=", "82×67="),
    @("91×29=", "62×14="),
    @("12×91=", "81×52="),
    @("15×35=", "89×23="),
    @("60×69=", "71×28="),
    @("25×22=", "20×52="),
    @("74×11=", "19×73="),
    @("88×30=", "91×94="),
    @("25×17=", "11×66="),
    @("49×67=", "57×45="),
    @("34×31=", "72×18="),
    @("88×69=", "65×62="),
    @("81×78=", "65×31="),
    @("33×45=", "15×29="),
    @("89×90=", "27×53="),
    @("16×48=", "84×15=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
